# Weekly price update: a new weekly record is inserted at row 10
# (pushing the existing rows 10-28 down to 11-29), matching the
# "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10; this shifts every
# row from 10..28 down by one (to 11..29) and extends the used range
# to A1:R29, keeping the inherited number formatting (date style) on
# column D.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with this week's record.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44571
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112052
$ws.Range("G10").Value = "Albahaca"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 950
$ws.Range("N10").Value = "$/paquete"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 950
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
